$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename subsector labels in column B (keep ids in column A and codes in column C unchanged)
$ws.Range("B2").Value  = "Agriculture"
$ws.Range("B3").Value  = "Manufacturing"
$ws.Range("B4").Value  = "Construction"
$ws.Range("B5").Value  = "Retail & wholesale"
$ws.Range("B6").Value  = "Transporting & storage"
$ws.Range("B7").Value  = "Hotel & restaurant"
$ws.Range("B8").Value  = "Information & communication"
$ws.Range("B9").Value  = "Finance"
$ws.Range("B10").Value = "Real estate and housing"
$ws.Range("B11").Value = "Professional tech. & science"
$ws.Range("B12").Value = "Administrative"
$ws.Range("B13").Value = "Public administration"
$ws.Range("B14").Value = "Education"
$ws.Range("B15").Value = "Health and social services"
$ws.Range("B16").Value = "Art, entertainment and recreation"
$ws.Range("B17").Value = "Other tertiary sectors"
$ws.Range("B18").Value = "Residential"

# Update the active selection to B18 (matches the saved view state in the target workbook)
$ws.Range("B18").Select()
